$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.765646815299988
$ws.Range("B1").Value = 2.100370645523071
$ws.Range("C1").Value = 2.250666618347168
$ws.Range("D1").Value = 2.966165781021118
$ws.Range("E1").Value = 1.699885964393616
